$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove formulas from C2 and C3, keep them as static values
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 67

# Row 7 rearrangement: A7 = "sdf" (was B7 "reg"), B7 = 1243 (new number), C7 stays "java.rmi.UnmarshalException"
$ws.Range("A7").Value = "sdf"
$ws.Range("B7").Value = 1243
$ws.Range("C7").Value = "java.rmi.UnmarshalException"

# New row 8
$ws.Range("A8").Value = 345
$ws.Range("B8").Value = "srg"
$ws.Range("C8").Value = "java.rmi.UnmarshalException"

# New row 9, just C9 with style s=1 (numFmt 0.0), empty value
$ws.Range("C9").NumberFormat = $ws.Range("C2").NumberFormat

# Update selection to C4
$ws.Range("C4").Select()
